# Update the answer cells in the multiplication-practice table.
# The table has data rows at 1, 5, 10, 15, 20 (blank spacer rows between),
# each with 5 columns of "A×B=C" answers. Cells are addressed by
# (row, column) rather than a blind text Find/Replace because some of the
# new values collide with old values located earlier in the document
# (e.g. row 1 col 4's new text equals row 10 col 1's old text).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "327×2=654"
$t.Cell(1, 2).Range.Text = "470×9=4230"
$t.Cell(1, 3).Range.Text = "778×6=4668"
$t.Cell(1, 4).Range.Text = "440×6=2640"
$t.Cell(1, 5).Range.Text = "693×5=3465"
$t.Cell(5, 1).Range.Text = "147×8=1176"
$t.Cell(5, 2).Range.Text = "493×6=2958"
$t.Cell(5, 3).Range.Text = "836×4=3344"
$t.Cell(5, 4).Range.Text = "110×2=220"
$t.Cell(5, 5).Range.Text = "693×2=1386"
$t.Cell(10, 1).Range.Text = "910×8=7280"
$t.Cell(10, 2).Range.Text = "119×8=952"
$t.Cell(10, 3).Range.Text = "681×2=1362"
$t.Cell(10, 4).Range.Text = "717×4=2868"
$t.Cell(10, 5).Range.Text = "725×5=3625"
$t.Cell(15, 1).Range.Text = "794×9=7146"
$t.Cell(15, 2).Range.Text = "924×2=1848"
$t.Cell(15, 3).Range.Text = "213×3=639"
$t.Cell(15, 4).Range.Text = "812×8=6496"
$t.Cell(15, 5).Range.Text = "948×4=3792"
$t.Cell(20, 1).Range.Text = "837×4=3348"
$t.Cell(20, 2).Range.Text = "171×9=1539"
$t.Cell(20, 3).Range.Text = "404×4=1616"
$t.Cell(20, 4).Range.Text = "207×8=1656"
$t.Cell(20, 5).Range.Text = "972×7=6804"
